$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric columns (I:T) for rows 2-7 with the newly computed TPM values.
$ws.Range("I2").Value = 0.6974138631561839
$ws.Range("J2").Value = 0.697413863156184
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.593519
$ws.Range("N2").Value = 43.780557
$ws.Range("O2").Value = 0.2780626807805203
$ws.Range("P2").Value = 0.2780626807805203
$ws.Range("Q2").Value = 0.6679842806806666
$ws.Range("R2").Value = 6.011858526126
$ws.Range("S2").Value = 0.1939247684027074
$ws.Range("T2").Value = 0.1939247684027075

$ws.Range("I3").Value = 0.6974138631561839
$ws.Range("J3").Value = 0.697413863156184
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.06491966666666
$ws.Range("N3").Value = 57.194759
$ws.Range("O3").Value = 0.3632600657441564
$ws.Range("P3").Value = 0.3632600657441565
$ws.Range("Q3").Value = 0.872652212929111
$ws.Range("R3").Value = 7.853869916361999
$ws.Range("S3").Value = 0.2533426057810015
$ws.Range("T3").Value = 0.2533426057810015

$ws.Range("I4").Value = 0.6974138631561839
$ws.Range("J4").Value = 0.697413863156184
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.82440066666667
$ws.Range("N4").Value = 56.473202
$ws.Range("O4").Value = 0.3586772534753233
$ws.Range("P4").Value = 0.3586772534753233
$ws.Range("Q4").Value = 0.861643016915111
$ws.Range("R4").Value = 7.754787152236
$ws.Range("S4").Value = 0.250146488972475
$ws.Range("T4").Value = 0.250146488972475

$ws.Range("I5").Value = 0.302586136843816
$ws.Range("J5").Value = 0.302586136843816
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 14.593519
$ws.Range("N5").Value = 43.780557
$ws.Range("O5").Value = 0.2780626807805203
$ws.Range("P5").Value = 0.2780626807805203
$ws.Range("Q5").Value = 0.2898175583273334
$ws.Range("R5").Value = 2.608358024946
$ws.Range("S5").Value = 0.08413791237781286
$ws.Range("T5").Value = 0.08413791237781286

$ws.Range("I6").Value = 0.302586136843816
$ws.Range("J6").Value = 0.302586136843816
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.06491966666666
$ws.Range("N6").Value = 57.194759
$ws.Range("O6").Value = 0.3632600657441564
$ws.Range("P6").Value = 0.3632600657441565
$ws.Range("Q6").Value = 0.3786165946335555
$ws.Range("R6").Value = 3.407549351702
$ws.Range("S6").Value = 0.1099174599631549
$ws.Range("T6").Value = 0.1099174599631549

$ws.Range("I7").Value = 0.302586136843816
$ws.Range("J7").Value = 0.302586136843816
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.82440066666667
$ws.Range("N7").Value = 56.473202
$ws.Range("O7").Value = 0.3586772534753233
$ws.Range("P7").Value = 0.3586772534753233
$ws.Range("Q7").Value = 0.3738400476395555
$ws.Range("R7").Value = 3.364560428756
$ws.Range("S7").Value = 0.1085307645028483
$ws.Range("T7").Value = 0.1085307645028483

# The MuSCs-as-sender rows (old rows 8-10) are no longer part of the
# recomputed TPM output, so remove them entirely - this also shrinks the
# sheet dimension from A1:T10 down to A1:T7.
$ws.Rows("8:10").Delete()
